$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.265.92'
$ws.Range('E2').Value = '  +0.87%  '

$ws.Range('D3').Value = '3.500.66'
$ws.Range('E3').Value = '  +0.36%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '''586.00'
$ws.Range('E5').Value = '  +0.41%  '

$ws.Range('D6').Value = '''134.03'
$ws.Range('E6').Value = '  +2.01%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  +0.23%  '

$ws.Range('E9').Value = '  +1.42%  '

$ws.Range('E10').Value = '  +1.80%  '

$ws.Range('E11').Value = '  +1.10%  '

$ws.Range('D12').Value = '4.102.27'
$ws.Range('E12').Value = '  +0.64%  '

$ws.Range('E13').Value = '  +1.03%  '

$ws.Range('D14').Value = '''0.0000181'
$ws.Range('E14').Value = '  +2.51%  '

$ws.Range('D15').Value = '3.502.82'
$ws.Range('E15').Value = '  +0.99%  '

$ws.Range('D16').Value = '''26.05'
$ws.Range('E16').Value = '  -5.25%  '

$ws.Range('D17').Value = '64.308.08'
$ws.Range('E17').Value = '  +0.93%  '

$ws.Range('E18').Value = '  -0.68%  '

$ws.Range('E19').Value = '  +1.27%  '

$ws.Range('D20').Value = '''13.71'
$ws.Range('E20').Value = '  -4.35%  '

$ws.Range('D21').Value = '''394.02'
$ws.Range('E21').Value = '  +2.60%  '

$ws.Range('D22').Value = '''0.572'
$ws.Range('E22').Value = '  -1.21%  '

$ws.Range('D23').Value = '3.641.59'
$ws.Range('E23').Value = '  +0.49%  '

$ws.Range('D24').Value = '''74.09'
$ws.Range('E24').Value = '  +1.48%  '

$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('D26').Value = '''5.68'
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').Value = '''0.0000114'
$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').Value = '''7.44'
$ws.Range('E28').Value = '  -1.34%  '

$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('D30').Value = '''1.52'
$ws.Range('E30').Value = '  -4.41%  '

$ws.Range('D31').Value = '''8.29'
$ws.Range('E31').Value = '  -0.24%  '

$ws.Range('D32').Value = '''2.22'
$ws.Range('E32').Value = '  -0.41%  '

$ws.Range('D33').Value = '3.521.16'
$ws.Range('E33').Value = '  +0.83%  '

$ws.Range('E34').Value = '  +0.02%  '

$ws.Range('E35').Value = '  +3.28%  '

$ws.Range('D36').Value = '''23.47'
$ws.Range('E36').Value = '  -0.23%  '

$ws.Range('D37').Value = '''5.20'
$ws.Range('E37').Value = '  -2.88%  '

$ws.Range('E38').Value = '  +0.00%  '

$ws.Range('D39').Value = '''6.92'
$ws.Range('E39').Value = '  -0.86%  '

$ws.Range('D40').Value = '''161.94'
$ws.Range('E40').Value = '  +0.28%  '

$ws.Range('D41').Value = '''0.0783'
$ws.Range('E41').Value = '  -1.91%  '

$ws.Range('E42').Value = '  -0.39%  '

$ws.Range('E43').Value = '  +0.20%  '

$ws.Range('D44').Value = '''25.17'
$ws.Range('E44').Value = '  -4.41%  '

$ws.Range('D45').Value = '''4.40'
$ws.Range('E45').Value = '  +0.05%  '

$ws.Range('D46').Value = '''1.17'
$ws.Range('E46').Value = '  -3.31%  '

$ws.Range('D47').Value = '''1.65'
$ws.Range('E47').Value = '  +1.57%  '

$ws.Range('D48').Value = '2.475.60'
$ws.Range('E48').Value = '  +2.46%  '

$ws.Range('E49').Value = '  -0.77%  '

$ws.Range('D50').Value = '''0.894'
$ws.Range('E50').Value = '  -0.60%  '

$ws.Range('E51').Value = '  -1.50%  '
